# Update the "想去人数" (column F) counts on both the "展览" and "全部类型"
# sheets, which hold duplicate copies of the same event listing.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1402
    11 = 4698
    18 = 4166
    19 = 936
    26 = 384
    30 = 50
    34 = 555
    37 = 6
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
